# Updates the cryptos list (rows 2-51) with the latest scraped price/volume
# figures. A couple of coins (InjectiveProtocol/Hedera and
# NEARProtocol/LidoDAOToken) also swapped rank positions, so their
# name/link/price/volume cells are rewritten wholesale for those rows.
#
# Price cells (column D) are stored as plain text in the source data (not
# numbers) - e.g. "109.80" must keep its trailing zero and "45.280.93" isn't
# a valid number at all. Whenever the new price string would otherwise be
# auto-coerced to a number by Excel, the cell is pre-formatted as Text ("@")
# so the literal string is preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '45.280.93'
$ws.Cells.Item(2, 5).Value = '  +5.14%  '

$ws.Cells.Item(3, 4).Value = '2.360.73'
$ws.Cells.Item(3, 5).Value = '  +2.10%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '109.80'
$ws.Cells.Item(5, 5).Value = '  +2.02%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '309.33'
$ws.Cells.Item(6, 5).Value = '  -0.85%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.629'
$ws.Cells.Item(7, 5).Value = '  +0.25%  '

$ws.Cells.Item(8, 5).Value = '  -0.36%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.617'
$ws.Cells.Item(9, 5).Value = '  +1.60%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '41.34'
$ws.Cells.Item(10, 5).Value = '  +2.44%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0918'
$ws.Cells.Item(11, 5).Value = '  +0.59%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '8.47'
$ws.Cells.Item(12, 5).Value = '  +1.14%  '

$ws.Cells.Item(13, 5).Value = '  +1.41%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '0.985'
$ws.Cells.Item(14, 5).Value = '  -1.06%  '

$ws.Cells.Item(15, 4).Value = '2.713.40'
$ws.Cells.Item(15, 5).Value = '  +1.82%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '15.36'
$ws.Cells.Item(16, 5).Value = '  -0.06%  '

$ws.Cells.Item(17, 4).Value = '2.346.38'
$ws.Cells.Item(17, 5).Value = '  +1.81%  '

$ws.Cells.Item(18, 4).Value = '45.152.49'
$ws.Cells.Item(18, 5).Value = '  +5.18%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '7.30'
$ws.Cells.Item(19, 5).Value = '  -2.34%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.0000107'
$ws.Cells.Item(20, 5).Value = '  +1.16%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.49'
$ws.Cells.Item(21, 5).Value = '  +2.79%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '73.32'
$ws.Cells.Item(22, 5).Value = '  -0.39%  '

$ws.Cells.Item(23, 5).Value = '  -1.84%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '259.07'
$ws.Cells.Item(24, 5).Value = '  -2.54%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.31'
$ws.Cells.Item(25, 5).Value = '  +2.53%  '

$ws.Cells.Item(26, 5).Value = '  -0.37%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.10'
$ws.Cells.Item(27, 5).Value = '  +0.85%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.38'
$ws.Cells.Item(28, 5).Value = '  -3.83%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.36'
$ws.Cells.Item(29, 5).Value = '  +2.71%  '

$ws.Cells.Item(30, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '38.40'
$ws.Cells.Item(30, 5).Value = '  -0.92%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '22.43'
$ws.Cells.Item(31, 5).Value = '  +0.26%  '

$ws.Cells.Item(32, 2).Value = 'Hedera'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.0961'
$ws.Cells.Item(32, 5).Value = '  +10.39%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '170.61'
$ws.Cells.Item(33, 5).Value = '  +2.85%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.93'
$ws.Cells.Item(34, 5).Value = '  +5.94%  '

$ws.Cells.Item(35, 5).Value = '  +0.64%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '4.84'
$ws.Cells.Item(36, 5).Value = '  +4.49%  '

$ws.Cells.Item(37, 5).Value = '  +1.42%  '

$ws.Cells.Item(38, 2).Value = 'NEARProtocol'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.93'
$ws.Cells.Item(38, 5).Value = '  +7.25%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.0357'
$ws.Cells.Item(39, 5).Value = '  +0.14%  '

$ws.Cells.Item(40, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.94'
$ws.Cells.Item(40, 5).Value = '  +4.11%  '

$ws.Cells.Item(41, 5).Value = '  +7.80%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '99.43'
$ws.Cells.Item(42, 5).Value = '  -5.13%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.232'
$ws.Cells.Item(43, 5).Value = '  +0.06%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '69.90'
$ws.Cells.Item(44, 5).Value = '  -1.40%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.999'
$ws.Cells.Item(45, 5).Value = '  -0.53%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '12.75'
$ws.Cells.Item(46, 5).Value = '  +2.98%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '81.86'
$ws.Cells.Item(47, 5).Value = '  +7.26%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '112.64'
$ws.Cells.Item(48, 5).Value = '  -0.02%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '9.29'
$ws.Cells.Item(49, 5).Value = '  +4.84%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '5.51'
$ws.Cells.Item(50, 5).Value = '  +4.80%  '

$ws.Cells.Item(51, 4).Value = '1.641.68'
$ws.Cells.Item(51, 5).Value = '  -1.83%  '
